{"js": "// Update the multiplication-fact table: 5 populated rows (0, 4, 9, 14, 19),\n// each with 5 cells (columns 0-4). Every populated cell's equation text is\n// replaced with a new equation, per the target diff. Cells are addressed by\n// (row, col) position in the table so that values which happen to collide\n// with another cell's old/new text (e.g. \"634\u00d79=5706\" moves from row0/col1\n// to row0/col3) are never mixed up.\nconst updates = [\n  { row: 0, col: 0, from: \"553\u00d76=3318\", to: \"936\u00d75=4680\" },\n  { row: 0, col: 1, from: \"634\u00d79=5706\", to: \"388\u00d76=2328\" },\n  { row: 0, col: 2, from: \"845\u00d74=3380\", to: \"347\u00d74=1388\" },\n  { row: 0, col: 3, from: \"395\u00d73=1185\", to: \"634\u00d79=5706\" },\n  { row: 0, col: 4, from: \"119\u00d78=952\", to: \"957\u00d77=6699\" },\n\n  { row: 4, col: 0, from: \"736\u00d77=5152\", to: \"739\u00d79=6651\" },\n  { row: 4, col: 1, from: \"220\u00d79=1980\", to: \"546\u00d72=1092\" },\n  { row: 4, col: 2, from: \"508\u00d77=3556\", to: \"915\u00d77=6405\" },\n  { row: 4, col: 3, from: \"582\u00d73=1746\", to: \"224\u00d76=1344\" },\n  { row: 4, col: 4, from: \"699\u00d72=1398\", to: \"947\u00d77=6629\" },\n\n  { row: 9, col: 0, from: \"662\u00d76=3972\", to: \"836\u00d75=4180\" },\n  { row: 9, col: 1, from: \"960\u00d74=3840\", to: \"784\u00d72=1568\" },\n  { row: 9, col: 2, from: \"436\u00d75=2180\", to: \"703\u00d74=2812\" },\n  { row: 9, col: 3, from: \"876\u00d77=6132\", to: \"367\u00d77=2569\" },\n  { row: 9, col: 4, from: \"508\u00d76=3048\", to: \"229\u00d73=687\" },\n\n  { row: 14, col: 0, from: \"308\u00d78=2464\", to: \"483\u00d72=966\" },\n  { row: 14, col: 1, from: \"529\u00d79=4761\", to: \"884\u00d75=4420\" },\n  { row: 14, col: 2, from: \"624\u00d75=3120\", to: \"263\u00d73=789\" },\n  { row: 14, col: 3, from: \"303\u00d73=909\", to: \"428\u00d79=3852\" },\n  { row: 14, col: 4, from: \"169\u00d74=676\", to: \"413\u00d79=3717\" },\n\n  { row: 19, col: 0, from: \"348\u00d77=2436\", to: \"686\u00d74=2744\" },\n  { row: 19, col: 1, from: \"799\u00d75=3995\", to: \"575\u00d74=2300\" },\n  { row: 19, col: 2, from: \"246\u00d72=492\", to: \"697\u00d72=1394\" },\n  { row: 19, col: 3, from: \"859\u00d73=2577\", to: \"382\u00d78=3056\" },\n  { row: 19, col: 4, from: \"793\u00d72=1586\", to: \"567\u00d76=3402\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nfor (const u of updates) {\n  const cell = table.getCell(u.row, u.col);\n  const results = cell.body.search(u.from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${u.from}\" in cell (${u.row},${u.col}), found ${results.items.length}`\n    );\n  }\n  results.items[0].insertText(u.to, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the multiplication-fact table: 5 populated rows (1, 5, 10, 15, 20 in\n# 1-based Word COM row numbering), each with 5 cells (columns 1-5). Every\n# populated cell's equation text is replaced with a new equation, per the\n# target diff. Cells are addressed by (row, col) table position so that\n# values which happen to collide with another cell's old/new text (e.g.\n# \"634\u00d79=5706\" moves from row1/col2 to row1/col4) are never mixed up.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @(\n  @{ Row = 1;  Col = 1; From = \"553\u00d76=3318\"; Text = \"936\u00d75=4680\" },\n  @{ Row = 1;  Col = 2; From = \"634\u00d79=5706\"; Text = \"388\u00d76=2328\" },\n  @{ Row = 1;  Col = 3; From = \"845\u00d74=3380\"; Text = \"347\u00d74=1388\" },\n  @{ Row = 1;  Col = 4; From = \"395\u00d73=1185\"; Text = \"634\u00d79=5706\" },\n  @{ Row = 1;  Col = 5; From = \"119\u00d78=952\";  Text = \"957\u00d77=6699\" },\n\n  @{ Row = 5;  Col = 1; From = \"736\u00d77=5152\"; Text = \"739\u00d79=6651\" },\n  @{ Row = 5;  Col = 2; From = \"220\u00d79=1980\"; Text = \"546\u00d72=1092\" },\n  @{ Row = 5;  Col = 3; From = \"508\u00d77=3556\"; Text = \"915\u00d77=6405\" },\n  @{ Row = 5;  Col = 4; From = \"582\u00d73=1746\"; Text = \"224\u00d76=1344\" },\n  @{ Row = 5;  Col = 5; From = \"699\u00d72=1398\"; Text = \"947\u00d77=6629\" },\n\n  @{ Row = 10; Col = 1; From = \"662\u00d76=3972\"; Text = \"836\u00d75=4180\" },\n  @{ Row = 10; Col = 2; From = \"960\u00d74=3840\"; Text = \"784\u00d72=1568\" },\n  @{ Row = 10; Col = 3; From = \"436\u00d75=2180\"; Text = \"703\u00d74=2812\" },\n  @{ Row = 10; Col = 4; From = \"876\u00d77=6132\"; Text = \"367\u00d77=2569\" },\n  @{ Row = 10; Col = 5; From = \"508\u00d76=3048\"; Text = \"229\u00d73=687\" },\n\n  @{ Row = 15; Col = 1; From = \"308\u00d78=2464\"; Text = \"483\u00d72=966\" },\n  @{ Row = 15; Col = 2; From = \"529\u00d79=4761\"; Text = \"884\u00d75=4420\" },\n  @{ Row = 15; Col = 3; From = \"624\u00d75=3120\"; Text = \"263\u00d73=789\" },\n  @{ Row = 15; Col = 4; From = \"303\u00d73=909\";  Text = \"428\u00d79=3852\" },\n  @{ Row = 15; Col = 5; From = \"169\u00d74=676\";  Text = \"413\u00d79=3717\" },\n\n  @{ Row = 20; Col = 1; From = \"348\u00d77=2436\"; Text = \"686\u00d74=2744\" },\n  @{ Row = 20; Col = 2; From = \"799\u00d75=3995\"; Text = \"575\u00d74=2300\" },\n  @{ Row = 20; Col = 3; From = \"246\u00d72=492\";  Text = \"697\u00d72=1394\" },\n  @{ Row = 20; Col = 4; From = \"859\u00d73=2577\"; Text = \"382\u00d78=3056\" },\n  @{ Row = 20; Col = 5; From = \"793\u00d72=1586\"; Text = \"567\u00d76=3402\" }\n)\n\nforeach ($u in $updates) {\n  $cellRange = $t.Cell($u.Row, $u.Col).Range\n  # Trim the trailing cell-mark / paragraph-mark (wdCharacter = 1) off the\n  # end of the cell range so only the visible run text is replaced; this\n  # keeps the original run/paragraph formatting (font, size, alignment)\n  # untouched instead of clearing it.\n  [void]$cellRange.MoveEnd(1, -1)\n  if ($cellRange.Text -ne $u.From) {\n    throw \"Cell ($($u.Row),$($u.Col)) expected '$($u.From)' but found '$($cellRange.Text)'\"\n  }\n  $cellRange.Text = $u.Text\n}\n"}
